$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gdf9"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3517556666666666
$ws.Range("H2").Value = 1.055267
$ws.Range("I2").Value = 0.04271818225758384
$ws.Range("J2").Value = 0.04271818225758384
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.707786666666667
$ws.Range("N2").Value = 5.12336
$ws.Range("O2").Value = 0.8764025646701329
$ws.Range("P2").Value = 0.8764025646701328
$ws.Range("Q2").Value = 0.6007236374577777
$ws.Range("R2").Value = 5.40651273712
$ws.Range("S2").Value = 0.03743832448859264
$ws.Range("T2").Value = 0.03743832448859263

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gdf9"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3517556666666666
$ws.Range("H3").Value = 1.055267
$ws.Range("I3").Value = 0.04271818225758384
$ws.Range("J3").Value = 0.04271818225758384
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.240846
$ws.Range("N3").Value = 0.7225379999999999
$ws.Range("O3").Value = 0.1235974353298672
$ws.Range("P3").Value = 0.1235974353298672
$ws.Range("Q3").Value = 0.08471894529399998
$ws.Range("R3").Value = 0.7624705076459999
$ws.Range("S3").Value = 0.005279857768991199
$ws.Range("T3").Value = 0.005279857768991198

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdf9"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.870075
$ws.Range("H4").Value = 20.610225
$ws.Range("I4").Value = 0.8343209329201149
$ws.Range("J4").Value = 0.8343209329201149
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.707786666666667
$ws.Range("N4").Value = 5.12336
$ws.Range("O4").Value = 0.8764025646701329
$ws.Range("P4").Value = 0.8764025646701328
$ws.Range("Q4").Value = 11.732622484
$ws.Range("R4").Value = 105.593602356
$ws.Range("S4").Value = 0.7312010053691665
$ws.Range("T4").Value = 0.7312010053691665

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdf9"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.870075
$ws.Range("H5").Value = 20.610225
$ws.Range("I5").Value = 0.8343209329201149
$ws.Range("J5").Value = 0.8343209329201149
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.240846
$ws.Range("N5").Value = 0.7225379999999999
$ws.Range("O5").Value = 0.1235974353298672
$ws.Range("P5").Value = 0.1235974353298672
$ws.Range("Q5").Value = 1.65463008345
$ws.Range("R5").Value = 14.89167075105
$ws.Range("S5").Value = 0.1031199275509484
$ws.Range("T5").Value = 0.1031199275509484

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Gdf9"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7357713333333334
$ws.Range("H6").Value = 2.207314
$ws.Range("I6").Value = 0.08935410825100797
$ws.Range("J6").Value = 0.08935410825100797
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.707786666666667
$ws.Range("N6").Value = 5.12336
$ws.Range("O6").Value = 0.8764025646701329
$ws.Range("P6").Value = 0.8764025646701328
$ws.Range("Q6").Value = 1.256540472782222
$ws.Range("R6").Value = 11.30886425504
$ws.Range("S6").Value = 0.07831016963499607
$ws.Range("T6").Value = 0.07831016963499605

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Gdf9"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7357713333333334
$ws.Range("H7").Value = 2.207314
$ws.Range("I7").Value = 0.08935410825100797
$ws.Range("J7").Value = 0.08935410825100797
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.240846
$ws.Range("N7").Value = 0.7225379999999999
$ws.Range("O7").Value = 0.1235974353298672
$ws.Range("P7").Value = 0.1235974353298672
$ws.Range("Q7").Value = 0.177207582548
$ws.Range("R7").Value = 1.594868242932
$ws.Range("S7").Value = 0.01104393861601191
$ws.Range("T7").Value = 0.01104393861601191

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gdf9"
$ws.Range("C8").Value = "Bmpr1b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2767293333333333
$ws.Range("H8").Value = 0.830188
$ws.Range("I8").Value = 0.03360677657129334
$ws.Range("J8").Value = 0.03360677657129334
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.707786666666667
$ws.Range("N8").Value = 5.12336
$ws.Range("O8").Value = 0.8764025646701329
$ws.Range("P8").Value = 0.8764025646701328
$ws.Range("Q8").Value = 0.4725946657422222
$ws.Range("R8").Value = 4.25335199168
$ws.Range("S8").Value = 0.02945306517737762
$ws.Range("T8").Value = 0.02945306517737762

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gdf9"
$ws.Range("C9").Value = "Bmpr1b"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2767293333333333
$ws.Range("H9").Value = 0.830188
$ws.Range("I9").Value = 0.03360677657129334
$ws.Range("J9").Value = 0.03360677657129334
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.240846
$ws.Range("N9").Value = 0.7225379999999999
$ws.Range("O9").Value = 0.1235974353298672
$ws.Range("P9").Value = 0.1235974353298672
$ws.Range("Q9").Value = 0.066649153016
$ws.Range("R9").Value = 0.5998423771439999
$ws.Range("S9").Value = 0.004153711393915725
$ws.Range("T9").Value = 0.004153711393915725
